$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.941.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.512.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.46%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("E8").Value = '  -1.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.514.52'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("E11").Value = '  +0.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.43'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.357'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.959.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.06'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.921.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.511.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.31%  '

$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.35%  '

$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.425'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.167'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.30%  '

$ws.Range("E27").Value = '  +0.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0768'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.44%  '

$ws.Range("E31").Value = '  -1.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '163.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.27%  '

$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("E34").Value = '  -4.77%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.42'
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.22'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.00%  '

$ws.Range("E39").Value = '  -0.18%  '

$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.800'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.83%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '278.52'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("E45").Value = '  +0.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.596'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("E47").Value = '  +0.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.20'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.76%  '

$ws.Range("E50").Value = '  -0.31%  '

$ws.Range("E51").Value = '  -1.86%  '
